# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the 10c728f3-... file on each
# locale sheet to reflect a freshly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-27 07:31:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-27 07:31:27"
